$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Bag. 3" (sheet4.xml): fill in Concatenate/Len and Left/Mid/Right
# formula columns that were previously blank.
# ---------------------------------------------------------------------------
$bag3 = $wb.Worksheets.Item("Bag. 3")

$bag3.Range("D2").Formula = "=CONCATENATE(C2,""-"",A2)"
$bag3.Range("E2").Formula = "=LEN(B2)"
$bag3.Range("D3").Formula = "=CONCATENATE(C3,""-"",A3)"
$bag3.Range("E3").Formula = "=LEN(B3)"
$bag3.Range("D4").Formula = "=CONCATENATE(C4,""-"",A4)"
$bag3.Range("E4").Formula = "=LEN(B4)"
$bag3.Range("D5").Formula = "=CONCATENATE(C5,""-"",A5)"
$bag3.Range("E5").Formula = "=LEN(B5)"

# The "Panjang Nama" (Len) column is centered, unlike the Concatenate column.
$bag3.Range("E2:E5").HorizontalAlignment = -4108
$bag3.Range("E2:E5").VerticalAlignment = -4108

$bag3.Range("C10").Formula = "=LEFT(B10,3)"
$bag3.Range("D10").Formula = "=MID(B10,4,7)"
$bag3.Range("E10").Formula = "=RIGHT(B10,2)"
$bag3.Range("C11").Formula = "=LEFT(B11,3)"
$bag3.Range("D11").Formula = "=MID(B11,4,7)"
$bag3.Range("E11").Formula = "=RIGHT(B11,2)"
$bag3.Range("C12").Formula = "=LEFT(B12,3)"
$bag3.Range("D12").Formula = "=MID(B12,4,7)"
$bag3.Range("E12").Formula = "=RIGHT(B12,2)"
$bag3.Range("C13").Formula = "=LEFT(B13,3)"
$bag3.Range("D13").Formula = "=MID(B13,4,7)"
$bag3.Range("E13").Formula = "=RIGHT(B13,2)"

# The Left/Mid/Right (NIK split) columns are centered too.
$bag3.Range("C10:E13").HorizontalAlignment = -4108
$bag3.Range("C10:E13").VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# Sheet "Bag. 5" (sheet6.xml): fill in Vlookup/Hlookup formula columns that
# were previously blank, referencing the "Data Referensi" lookup tables.
# ---------------------------------------------------------------------------
$bag5 = $wb.Worksheets.Item("Bag. 5")

for ($r = 2; $r -le 11; $r++) {
    $bag5.Range("C$r").Formula = "=VLOOKUP(B$r,'Data Referensi'!`$A`$2:`$D`$11,2,FALSE)"
    $bag5.Range("D$r").Formula = "=VLOOKUP(B$r,'Data Referensi'!`$A`$2:`$D`$11,4,FALSE)"
    $bag5.Range("E$r").Formula = "=HLOOKUP(D$r,'Data Referensi'!`$B`$14:`$E`$15,2,FALSE)"
}

# ---------------------------------------------------------------------------
# Restore/update each sheet's view state (zoom + selection) to match the
# author's last-saved session, then leave "Bag. 5" as the active tab.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Activate()
$summary.Range("F5").Select()
$excel.ActiveWindow.Zoom = 97

$bag1 = $wb.Worksheets.Item("Bag. 1")
$bag1.Activate()
$bag1.Range("C14").Select()
$excel.ActiveWindow.Zoom = 116

$bag2 = $wb.Worksheets.Item("Bag. 2")
$bag2.Activate()
$bag2.Range("C5:D5").Select()
$excel.ActiveWindow.Zoom = 85

$bag3.Activate()
$bag3.Range("F21").Select()
$excel.ActiveWindow.Zoom = 61

$bag4 = $wb.Worksheets.Item("Bag. 4")
$bag4.Activate()
$bag4.Range("O11").Select()
$excel.ActiveWindow.Zoom = 70

$dataRef = $wb.Worksheets.Item("Data Referensi")
$dataRef.Activate()
$dataRef.Range("A15").Select()
$excel.ActiveWindow.Zoom = 38

$bag5.Activate()
$bag5.Range("E2").Select()
$excel.ActiveWindow.Zoom = 108
